# Update CDA Logical model for ST.r2b
# - Bump Version and Date metadata values
# - Insert a new "Jurisdiction" row into the Metadata table (before "Description")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version (row 3) and Date (row 8) values
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row before the "Description" row (row 11) and populate it
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
